{"js": "// Locate the (only) table in the document body.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\ntable.rows.load(\"items\");\nawait context.sync();\n\n// Helper: replace the first (and in this document, only) occurrence of\n// `find` with `replacement`, touching just that text run.\nasync function replaceText(find, replacement) {\n  const results = context.document.body.search(find, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  results.items[0].insertText(replacement, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// Row \"Variation within populations\": row height 612 -> 614 twips\n// (30.6pt -> 30.7pt), label renamed, and the Std.Obs value updated.\ntable.rows.items[1].preferredHeight = 30.7;\nawait context.sync();\nawait replaceText(\"Variation within populations\", \"Variation within sampling sites\");\nawait replaceText(\"-22.924\", \"-23.015\");\n\n// Row \"Variation between populations\": row height 612 -> 614 twips\n// (30.6pt -> 30.7pt), label renamed, and the Std.Obs value updated.\ntable.rows.items[2].preferredHeight = 30.7;\nawait context.sync();\nawait replaceText(\"Variation between populations\", \"Variation between sampling sites\");\nawait replaceText(\"23.138\", \"22.859\");\n\n// Row \"Variation between urban/rural groups\": Std.Obs and p updated.\nawait replaceText(\"0.282\", \"0.391\");\nawait replaceText(\"0.376\", \"0.336\");\n", "ps1": "# There is a single table in the document.\n$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n\n# Row 2 (1-based): \"Variation within populations\"\n#  - row height 612 -> 614 twips (30.6pt -> 30.7pt)\n#  - label renamed to \"Variation within sampling sites\"\n#  - Std.Obs value -22.924 -> -23.015\n$row1 = $table.Rows.Item(2)\n$row1.Height = 30.7\n$table.Cell(2, 1).Range.Text = \"Variation within sampling sites\"\n$table.Cell(2, 3).Range.Text = \"-23.015\"\n\n# Row 3 (1-based): \"Variation between populations\"\n#  - row height 612 -> 614 twips (30.6pt -> 30.7pt)\n#  - label renamed to \"Variation between sampling sites\"\n#  - Std.Obs value 23.138 -> 22.859\n$row2 = $table.Rows.Item(3)\n$row2.Height = 30.7\n$table.Cell(3, 1).Range.Text = \"Variation between sampling sites\"\n$table.Cell(3, 3).Range.Text = \"22.859\"\n\n# Row 4 (1-based): \"Variation between urban/rural groups\"\n#  - Std.Obs value 0.282 -> 0.391\n#  - p value 0.376 -> 0.336\n$table.Cell(4, 3).Range.Text = \"0.391\"\n$table.Cell(4, 5).Range.Text = \"0.336\"\n"}
